$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the bold/border/centered formatting used by the other header cells
# (e.g. H1) by copying its format onto the new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data cells I2 and J2, both numeric value 9 (plain, unstyled like the
# other data cells in row 2)
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
